$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarter-end years represented in column A (rows 2-22) were stored as
# date serials with a custom "YYYY-MM-DD HH:MM:SS" number format. The
# naive forecaster bugfix switches them to plain text labels ("2004Q4",
# "2005Q4", ...), matching the header's plain (non-date) style.
$years = 2004..2024
$dataRange = $ws.Range("A2:A22")

# Replace the date values with text quarter labels first.
for ($i = 0; $i -lt $years.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "$($years[$i])Q4"
}

# Re-use the header cell's style (border/font/alignment, General number
# format) instead of the old date-formatted style, so no stray custom
# number format / style entries are left behind.
$ws.Range("A1").Copy()
$dataRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
